$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.451.67"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").Value = "3.425.89"
$ws.Range("E3").Value = "  -2.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.18%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.422.50"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.95%  "
$ws.Range("E10").Value = "  -7.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.31%  "
$ws.Range("E12").Value = "  -7.88%  "
$ws.Range("D13").Value = "4.000.53"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("E14").Value = "  -8.10%  "
$ws.Range("D15").Value = "3.430.42"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.31%  "
$ws.Range("D18").Value = "64.400.94"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -11.43%  "
$ws.Range("E20").Value = "  -7.91%  "
$ws.Range("E21").Value = "  -6.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "379.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.542"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.84%  "
$ws.Range("D27").Value = "3.560.78"
$ws.Range("E27").Value = "  -2.81%  "
$ws.Range("E28").Value = "  -7.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.80%  "
$ws.Range("E31").Value = "  -9.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.15%  "
$ws.Range("D33").Value = "3.434.44"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -5.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.140"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "171.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.79%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.58%  "
$ws.Range("E40").Value = "  -10.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0759"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.800"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.21%  "
$ws.Range("E46").Value = "  -13.99%  "
$ws.Range("E47").Value = "  -9.98%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.72%  "
$ws.Range("D51").Value = "2.195.83"
$ws.Range("E51").Value = "  -4.98%  "
